$wb = $excel.ActiveWorkbook

# ALC!H103:M103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 601
$ws.Range("I103").Value = 601
$ws.Range("K103").Value = 1803
$ws.Range("M103").Value = -1217

# ALC!H116:N116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2083.75
$ws.Range("I116").Value = 1968.3334
$ws.Range("J116").Value = 2122.2222
$ws.Range("K116").Value = 1968.3334
$ws.Range("L116").Value = 2122.2222
$ws.Range("M116").Value = 1473.6666
$ws.Range("N116").Value = -9006.2222

# ALC!H134:N134
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 31392.38
$ws.Range("J134").Value = 31392.38
$ws.Range("L134").Value = 31392.38
$ws.Range("N134").Value = -41532.38

# ALC!H136:N136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 28064.445
$ws.Range("J136").Value = 28064.445
$ws.Range("L136").Value = 28064.445
$ws.Range("N136").Value = -38264.445

# ALC!H137:N137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1739.5652
$ws.Range("I137").Value = 1149.8
$ws.Range("J137").Value = 2845.375
$ws.Range("K137").Value = 3449.4
$ws.Range("L137").Value = 8536.125
$ws.Range("M137").Value = -899.3999999999996
$ws.Range("N137").Value = -13636.125

# ALC!H138:N138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3420.758
$ws.Range("I138").Value = 1807.7826
$ws.Range("J138").Value = 4372
$ws.Range("K138").Value = 5423.3478
$ws.Range("L138").Value = 13116
$ws.Range("M138").Value = -283.3477999999996
$ws.Range("N138").Value = -23396

# ARM!H32:N32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4377
$ws.Range("I32").Value = 3998.7058
$ws.Range("J32").Value = 5511.8823
$ws.Range("K32").Value = 3998.7058
$ws.Range("L32").Value = 5511.8823
$ws.Range("M32").Value = -3711.7058
$ws.Range("N32").Value = -6085.8823

# ARM!H50:N50
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 1485.8
$ws.Range("I50").Value = 548
$ws.Range("J50").Value = 2111
$ws.Range("K50").Value = 548
$ws.Range("L50").Value = 2111
$ws.Range("M50").Value = 166
$ws.Range("N50").Value = -3539

# ARM!H74:N74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1800.1562
$ws.Range("I74").Value = 1518.9048
$ws.Range("J74").Value = 2337.0908
$ws.Range("K74").Value = 1518.9048
$ws.Range("L74").Value = 2337.0908
$ws.Range("M74").Value = -644.9048
$ws.Range("N74").Value = -4085.0908

# ARM!H77:N77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1800.1562
$ws.Range("I77").Value = 1518.9048
$ws.Range("J77").Value = 2337.0908
$ws.Range("K77").Value = 7594.524
$ws.Range("L77").Value = 11685.454
$ws.Range("M77").Value = -3226.524
$ws.Range("N77").Value = -20421.454

# ARM!H122:N122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1510300.6
$ws.Range("I122").Value = 2850623.2
$ws.Range("J122").Value = 2437.5
$ws.Range("K122").Value = 8551869.600000001
$ws.Range("L122").Value = 7312.5
$ws.Range("M122").Value = -8549419.600000001
$ws.Range("N122").Value = -12212.5

# ARM!H132:N132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2443.0527
$ws.Range("I132").Value = 2170.3635
$ws.Range("J132").Value = 4242.8
$ws.Range("K132").Value = 6511.0905
$ws.Range("L132").Value = 12728.4
$ws.Range("M132").Value = -3981.0905
$ws.Range("N132").Value = -17788.4

# BSM!H105:N105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13902.883
$ws.Range("I105").Value = 19652.637
$ws.Range("J105").Value = 3361.6667
$ws.Range("K105").Value = 19652.637
$ws.Range("L105").Value = 3361.6667
$ws.Range("M105").Value = -17905.637
$ws.Range("N105").Value = -6855.6667

# CRP!H4:N4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 25000
$ws.Range("I4").Value = 25000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 25000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -24888
$ws.Range("N4").ClearContents()

# CRP!H31:N31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3246.5193
$ws.Range("J31").Value = 3270.7058
$ws.Range("L31").Value = 3270.7058
$ws.Range("N31").Value = -3860.7058

# CRP!H34:N34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3246.5193
$ws.Range("J34").Value = 3270.7058
$ws.Range("L34").Value = 3270.7058
$ws.Range("N34").Value = -3674.7058

# CRP!H59:N59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 26244.111
$ws.Range("J59").Value = 26244.111
$ws.Range("L59").Value = 26244.111
$ws.Range("N59").Value = -28534.111

# CRP!H94:N94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1674.7046
$ws.Range("I94").Value = 1493.7778
$ws.Range("J94").Value = 1799.9615
$ws.Range("K94").Value = 1493.7778
$ws.Range("L94").Value = 1799.9615
$ws.Range("M94").Value = -1042.7778
$ws.Range("N94").Value = -2701.9615

# CUL!H63:N63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4060
$ws.Range("I63").Value = 1900
$ws.Range("J63").Value = 5500
$ws.Range("K63").Value = 5700
$ws.Range("L63").Value = 16500
$ws.Range("M63").Value = -4951
$ws.Range("N63").Value = -17998

# CUL!H66:N66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 4060
$ws.Range("I66").Value = 1900
$ws.Range("J66").Value = 5500
$ws.Range("K66").Value = 17100
$ws.Range("L66").Value = 49500
$ws.Range("M66").Value = -13356
$ws.Range("N66").Value = -56988

# CUL!H68:N68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4457.727
$ws.Range("J68").Value = 1458.2273
$ws.Range("L68").Value = 4374.6819
$ws.Range("N68").Value = -5996.6819

# CUL!H71:N71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 4457.727
$ws.Range("J71").Value = 1458.2273
$ws.Range("L71").Value = 13124.0457
$ws.Range("N71").Value = -21236.0457

# CUL!H100:N100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 4750
$ws.Range("J100").Value = 4750
$ws.Range("L100").Value = 14250
$ws.Range("N100").Value = -15872

# CUL!H113:N113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1304827.6
$ws.Range("I113").Value = 1613360.6
$ws.Range("J113").Value = 667192.8
$ws.Range("K113").Value = 4840081.800000001
$ws.Range("L113").Value = 2001578.4
$ws.Range("M113").Value = -4837911.800000001
$ws.Range("N113").Value = -2005918.4

# CUL!H122:N122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 362.5
$ws.Range("I122").Value = 362.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3262.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -812.5
$ws.Range("N122").ClearContents()

# CUL!H132:N132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13250.5
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 14667.223
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 132005.007
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -137065.007

# GSM!H107:N107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1511.04
$ws.Range("I107").Value = 1053.6842
$ws.Range("J107").Value = 2959.3333
$ws.Range("K107").Value = 1053.6842
$ws.Range("L107").Value = 2959.3333
$ws.Range("M107").Value = 866.3158000000001
$ws.Range("N107").Value = -6799.3333

# GSM!H113:N113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 55556644
$ws.Range("I113").Value = 83334330
$ws.Range("J113").Value = 1268.8334
$ws.Range("K113").Value = 83334330
$ws.Range("L113").Value = 1268.8334
$ws.Range("M113").Value = -83332160
$ws.Range("N113").Value = -5608.8334

# GSM!H132:N132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6738.4614
$ws.Range("I132").Value = 15333.333
$ws.Range("J132").Value = 4160
$ws.Range("K132").Value = 45999.999
$ws.Range("L132").Value = 12480
$ws.Range("M132").Value = -43469.999
$ws.Range("N132").Value = -17540
